# Insert a new weekly record as row 26 (pushing the existing rows 26-42 down to 27-43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 45062
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112044
$ws.Cells.Item(26, 7).Value = "Perejil"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 380
$ws.Cells.Item(26, 11).Value = 1800
$ws.Cells.Item(26, 12).Value = 2000
$ws.Cells.Item(26, 13).Value = 1895
$ws.Cells.Item(26, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 948
$ws.Cells.Item(26, 17).Value = 2
$ws.Cells.Item(26, 18).Value = "Hortaliza"

Write-Host "Inserted new row 26 and shifted subsequent rows down."
